$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "Sheet2" worksheet after the existing sheets
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet2"

# ---------------------------------------------------------------------
# 2. Populate the new sheet with its text content
# ---------------------------------------------------------------------
$ws3.Range("A2").Value = "Most Active Days"
$ws3.Range("A3").Value = "Most Inactive days"
$ws3.Range("A4").Value = "Average Calories burned per day"
$ws3.Range("B5").Value = "Caloeries outliers analysis"
$ws3.Range("A6").Value = "Healthy sleep days"
$ws3.Range("A7").Value = "Unhealthy sleep days"
$ws3.Range("A8").Value = "People who sleep more than once a day"
$ws3.Range("B9").Value = "Sleeping pattern of those people and split of sleep each time"
$ws3.Range("A10").Value = "Calories burned to weight lost analysis"
$ws3.Range("A12").Value = "Activity split per day"
$ws3.Range("A13").Value = "Most Active time of day"
$ws3.Range("A14").Value = "Most Sedentary time of"

# Column widths (characters)
$ws3.Columns.Item(1).ColumnWidth = 53.04
$ws3.Columns.Item(2).ColumnWidth = 80.92

# The sheet opens with A14 selected / active
$ws3.Range("A14").Select()

# ---------------------------------------------------------------------
# 3. Changelog sheet: selection moves from F8 to F7 (and loses the
#    "tabSelected" flag automatically once Sheet2 becomes active)
# ---------------------------------------------------------------------
$wsChangelog = $wb.Worksheets.Item("Changelog")
$wsChangelog.Activate()
$wsChangelog.Range("F7").Select()

# ---------------------------------------------------------------------
# 4. Make the newly added Sheet2 the active / displayed tab
# ---------------------------------------------------------------------
$ws3.Activate()
